$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.XValues = "Time"
$s1.Values = "Voltage"
$s1.Name = "Voltage"
